# Add a new notebook type: VCP prediction using feature prototypes
# (only vote prototypes existed before).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# --- New header cell + k= column headers for the new "PRTS" block ---
$ws1.Range("Q12").Value = "BASADO EN PRTS (Promedio de características)"

$ws1.Range("Q13").Value = "k=1"
$ws1.Range("R13").Value = "k=3"
$ws1.Range("S13").Value = "k=5"
$ws1.Range("T13").Value = "k=7"
$ws1.Range("U13").Value = "k=9"
$ws1.Range("V13").Value = "k=11"

# --- New data row (INCV3-EUCLID) for the PRTS block ---
$ws1.Range("P19").Value = "INCV3-EUCLID"
$ws1.Range("Q19").Value = 4.91
$ws1.Range("R19").Value = 4.55
$ws1.Range("S19").Value = 4.81
$ws1.Range("T19").Value = 4.72
$ws1.Range("U19").Value = 4.57
$ws1.Range("V19").Value = 4.55

# --- Add the second (currently empty) worksheet, placed after Sheet1 ---
$newSheet = $wb.Worksheets.Add($null, $ws1)
$newSheet.Name = "Hoja1"

# --- Restore Sheet1 as the active / visible sheet, and reproduce the
#     view state (scroll position + selection) recorded in the workbook ---
$ws1.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 10
$win.ScrollColumn = 13
$ws1.Range("Q18").Select() | Out-Null
